# Auto-generated edit script: updates the cryptos price table.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, preserving the cells original
# style (Excel would otherwise auto-coerce numeric-looking strings like
# "8.00" or "0.0550" into numbers and drop significant trailing zeros).
function Set-TextValue($rng, [string]$val) {
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "63.768.60"
Set-TextValue $ws.Range("E2") "  +1.23%  "

Set-TextValue $ws.Range("D3") "2.655.40"
Set-TextValue $ws.Range("E3") "  +2.69%  "

Set-TextValue $ws.Range("E4") "  -0.01%  "

Set-TextValue $ws.Range("D5") "594.16"

Set-TextValue $ws.Range("E6") "  -0.91%  "

Set-TextValue $ws.Range("D8") "0.591"
Set-TextValue $ws.Range("E8") "  -1.15%  "

Set-TextValue $ws.Range("D9") "0.108"
Set-TextValue $ws.Range("E9") "  +0.01%  "

Set-TextValue $ws.Range("E10") "  -0.48%  "

Set-TextValue $ws.Range("E11") "  -0.06%  "

Set-TextValue $ws.Range("E12") "  +0.58%  "

Set-TextValue $ws.Range("D13") "27.49"
Set-TextValue $ws.Range("E13") "  +0.65%  "

Set-TextValue $ws.Range("D14") "3.131.13"
Set-TextValue $ws.Range("E14") "  +2.70%  "

Set-TextValue $ws.Range("D15") "63.686.87"
Set-TextValue $ws.Range("E15") "  +1.32%  "

Set-TextValue $ws.Range("E16") "  -0.05%  "

Set-TextValue $ws.Range("D17") "2.648.74"
Set-TextValue $ws.Range("E17") "  +2.38%  "

Set-TextValue $ws.Range("D18") "11.39"
Set-TextValue $ws.Range("E18") "  +0.66%  "

Set-TextValue $ws.Range("D19") "342.62"
Set-TextValue $ws.Range("E19") "  -0.20%  "

Set-TextValue $ws.Range("D20") "4.36"
Set-TextValue $ws.Range("E20") "  -0.62%  "

Set-TextValue $ws.Range("D21") "6.78"
Set-TextValue $ws.Range("E21") "  +1.77%  "

Set-TextValue $ws.Range("E22") "  +0.07%  "

Set-TextValue $ws.Range("D23") "68.07"
Set-TextValue $ws.Range("E23") "  +1.09%  "

Set-TextValue $ws.Range("E24") "  +13.33%  "

Set-TextValue $ws.Range("D25") "1.68"
Set-TextValue $ws.Range("E25") "  +5.46%  "

Set-TextValue $ws.Range("D26") "575.08"
Set-TextValue $ws.Range("E26") "  +23.36%  "

Set-TextValue $ws.Range("E27") "  -1.31%  "

Set-TextValue $ws.Range("D28") "8.58"
Set-TextValue $ws.Range("E28") "  +3.01%  "

Set-TextValue $ws.Range("B29") "Aptos"
Set-TextValue $ws.Range("C29") "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
Set-TextValue $ws.Range("D29") "8.00"
Set-TextValue $ws.Range("E29") "  +2.02%  "

Set-TextValue $ws.Range("B30") "Binance-PegBSC-USD"
Set-TextValue $ws.Range("C30") "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
Set-TextValue $ws.Range("D30") "1.00"
Set-TextValue $ws.Range("E30") "  +0.21%  "

Set-TextValue $ws.Range("D31") "2.00"
Set-TextValue $ws.Range("E31") "  +3.14%  "

Set-TextValue $ws.Range("E32") "  +11.40%  "

Set-TextValue $ws.Range("D33") "0.0₃0816"
Set-TextValue $ws.Range("E33") "  -1.01%  "

Set-TextValue $ws.Range("D34") "175.27"
Set-TextValue $ws.Range("E34") "  +0.14%  "

Set-TextValue $ws.Range("E35") "  +0.04%  "

Set-TextValue $ws.Range("E36") "  +0.03%  "

Set-TextValue $ws.Range("B37") "NEARProtocol"
Set-TextValue $ws.Range("C37") "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
Set-TextValue $ws.Range("D37") "4.69"
Set-TextValue $ws.Range("E37") "  +3.11%  "

Set-TextValue $ws.Range("B38") "EthereumClassic"
Set-TextValue $ws.Range("C38") "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
Set-TextValue $ws.Range("D38") "19.13"
Set-TextValue $ws.Range("E38") "  -0.02%  "

Set-TextValue $ws.Range("E39") "  +2.84%  "

Set-TextValue $ws.Range("D40") "171.75"
Set-TextValue $ws.Range("E40") "  +8.47%  "

Set-TextValue $ws.Range("D41") "0.999"
Set-TextValue $ws.Range("E41") "  -0.01%  "

Set-TextValue $ws.Range("D42") "40.44"
Set-TextValue $ws.Range("E42") "  +2.63%  "

Set-TextValue $ws.Range("E43") "  -0.65%  "

Set-TextValue $ws.Range("D44") "21.81"
Set-TextValue $ws.Range("E44") "  +2.51%  "

Set-TextValue $ws.Range("D45") "0.628"
Set-TextValue $ws.Range("E45") "  -1.59%  "

Set-TextValue $ws.Range("D46") "0.0550"
Set-TextValue $ws.Range("E46") "  +1.07%  "

Set-TextValue $ws.Range("D47") "0.0960"
Set-TextValue $ws.Range("E47") "  -0.89%  "

Set-TextValue $ws.Range("E48") "  +0.84%  "

Set-TextValue $ws.Range("D49") "18.65"
Set-TextValue $ws.Range("E49") "  +1.09%  "

Set-TextValue $ws.Range("E50") "  +1.75%  "

Set-TextValue $ws.Range("B51") "BabyDogeCoin"
Set-TextValue $ws.Range("C51") "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-TextValue $ws.Range("D51") "0.0₆0219"
Set-TextValue $ws.Range("E51") "  +13.86%  "

